$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Rushing sheet: update week-over-week cumulative rushing totals.
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

$rushing.Range("C2").Value = 0
$rushing.Range("E2").Value = 0

$rushing.Range("E3").Value = 0
$rushing.Range("F3").Value = 0

$rushing.Range("C4").Value = 96
$rushing.Range("D4").Value = 60

$rushing.Range("D5").Value = 34

$rushing.Range("C6").Value = 5
$rushing.Range("D6").Value = 4
$rushing.Range("E6").Value = 3

$rushing.Range("C11").Value = 0
$rushing.Range("D11").Value = 0
$rushing.Range("E11").Value = 0
$rushing.Range("F11").Value = 0

$rushing.Range("F4").Select()

# ---------------------------------------------------------------------------
# Receiving sheet: update cumulative receiving totals for several players,
# then add a new Week 15 row for A.Janovich (his first targets of the year).
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("C2").Value = 12
$receiving.Range("D2").Value = 9

$receiving.Range("C3").Value = 26
$receiving.Range("D3").Value = 21

$receiving.Range("C4").Value = 5
$receiving.Range("D4").Value = 5

$receiving.Range("C7").Value = 0
$receiving.Range("D7").Value = 0
$receiving.Range("E7").Value = 0
$receiving.Range("F7").Value = 0
$receiving.Range("G7").Value = 0
$receiving.Range("H7").Value = 0

$receiving.Range("C8").Value = 20
$receiving.Range("D8").Value = 15
$receiving.Range("E8").Value = 15
$receiving.Range("F8").Value = 8

$receiving.Range("C12").Value = 0
$receiving.Range("D12").Value = 0
$receiving.Range("E12").Value = 0
$receiving.Range("F12").Value = 0
$receiving.Range("G12").Value = 0
$receiving.Range("H12").Value = 0

# Insert a new row for A.Janovich between D.Felton (row 5) and J.Stanton (row 6)
# by shifting rows 6..14 down to 7..15 (bottom-up, so nothing is clobbered).
for ($r = 14; $r -ge 6; $r--) {
    $src = $receiving.Range("A" + $r + ":H" + $r)
    $dst = $receiving.Range("A" + ($r + 1) + ":H" + ($r + 1))
    $src.Copy($dst)
}
$excel.CutCopyMode = $false

$receiving.Range("A6").Value = 4
$receiving.Range("B6").Value = "A.Janovich"
$receiving.Range("C6").Value = 1
$receiving.Range("D6").Value = 1
$receiving.Range("E6").Value = 0
$receiving.Range("F6").Value = 0
$receiving.Range("G6").Value = 0
$receiving.Range("H6").Value = 0

# Renumber the player index column for every row pushed down by the insert.
$receiving.Range("A7").Value = 5
$receiving.Range("A8").Value = 6
$receiving.Range("A9").Value = 7
$receiving.Range("A10").Value = 8
$receiving.Range("A11").Value = 9
$receiving.Range("A12").Value = 10
$receiving.Range("A13").Value = 11
$receiving.Range("A14").Value = 12
$receiving.Range("A15").Value = 13

$receiving.Range("H15").Select()
$receiving.Activate()
